# "added TCS Food list modal"
# Appends the new TCS-food-category key/value rows (Meats, Poultry, Seafood,
# Bakery, Dairy, Pasta, Eggs, Fruits/Vegetables, Misc.) used to populate the
# new "TCS Food list" modal to the language sheet, directly below the
# existing key/value rows (which end at row 28).
#
# NOTE: the source strings use a literal two-character "\n" sequence (not an
# actual newline) as their line-break marker, matching the rest of the sheet
# (e.g. the existing "tcs_foods" value). PowerShell's backtick-n is a true
# newline, so plain backslash-n is used here instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("tcsFoodDesc_meats_title", "Meats"),
    @("tcsFoodDesc_meats", "Bacon - in raw form\nBeef - ground, roasts, steak\nGravy\nGround meats - all\nHot Dogs\nLunch meat\nMeat casseroles\nPork - ground, ham, roasts\nProcessed meats - all\nSausage\nSoups Stews"),
    @("tcsFoodDesc_poultry_title", "Poultry"),
    @("tcsFoodDesc_poultry", "Chicken - ground, roasted, barbequed, fried, nuggets, patties, strips\nCasseroles with chicken/turkey\nDressing\nGravy\nPrecooked, processed products\nTurkey – ground, roast\nSoups Stews"),
    @("tcsFoodDesc_seafood_title", "Seafood"),
    @("tcsFoodDesc_seafood", "Fish, Salmon, Tuna"),
    @("tcsFoodDesc_bakery_title", "Bakery Foods"),
    @("tcsFoodDesc_bakery", "Cream pastries\nCream/custard pies and tarts\nPudding - prepared"),
    @("tcsFoodDesc_dairy_title", "Dairy Foods"),
    @("tcsFoodDesc_dairy", "Whipped butter/whipped margarine\nCheese - mozzarella, cottage, cream, Ricotta\nCream - real, sauce, white\nDairy whipped topping\nIce cream\nMilk"),
    @("tcsFoodDesc_pasta_title", "Pasta"),
    @("tcsFoodDesc_pasta", "Noodles - all kinds, cooked\nRice - cooked"),
    @("tcsFoodDesc_eggs_title", "Eggs"),
    @("tcsFoodDesc_eggs", "Egg casseroles\nEgg dishes\nDeviled eggs\nFried eggs\nHard-cooked eggs\nOmelets\nScrambled eggs"),
    @("tcsFoodDesc_fruits_veg_title", "Fruits and Vegetables"),
    @("tcsFoodDesc_fruits_veg", "Dry beans - cooked\nPotatoes - baked, boiled mashed (fresh, instant), scalloped/augratin (fresh, dehydrated)\nCut/prepared fresh fruits and vegetables (including melons, tomatoes and salad greens)"),
    @("tcsFoodDesc_misc_title", "Misc."),
    @("tcsFoodDesc_misc", "Salad dressings prepared from a mix")
)

$startRow = 29
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $rows[$i][0]
    $ws.Range("B$r").Value = $rows[$i][1]
}

# Matches the author's final on-screen selection/scroll position when the
# workbook was saved (row 44, with the view scrolled down so row 28 is at
# the top) rather than the very last row that was written.
$ws.Range("B44").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
